$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old "profit" column (D) entirely - the sheet shrinks from 4 to 3 columns.
$ws.Range("A1:D5").ClearContents()

# New header row
$ws.Range("A1").Value = "school"
$ws.Range("B1").Value = "address"
$ws.Range("C1").Value = "students"

# New data rows
$ws.Range("A2").Value = "Portland State University"
$ws.Range("B2").Value = "1825 SW Broadway, Portland, OR 97201"
$ws.Range("C2").Value = 12490

$ws.Range("A3").Value = "University of Portland"
$ws.Range("B3").Value = "5000 N Willamette Blvd, Portland, OR 97203"
$ws.Range("C3").Value = 3700

$ws.Range("A4").Value = "Reed College"
$ws.Range("B4").Value = "3203 SE Woodstock Blvd, Portland, OR 97202"
$ws.Range("C4").Value = 1458

$ws.Range("A5").Value = "Lewis and Clark College"
$ws.Range("B5").Value = "615 S Palatine Hill Rd, OR 97219"
$ws.Range("C5").Value = 3520

# Column D no longer used - delete it
$ws.Columns.Item(4).Delete()

# Adjust column widths: B needs to widen to fit the new longest address
$ws.Columns.Item(2).ColumnWidth = 38

$ws.Range("E11").Select()
